$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = '[Na%Zhu%NULL%0,               Dingyu%Zhang%NULL%1,               Wenling%Wang%NULL%1,               Xingwang%Li%NULL%0,               Bo%Yang%NULL%3,               Jingdong%Song%NULL%1,               Xiang%Zhao%NULL%1,               Baoying%Huang%NULL%1,               Weifeng%Shi%NULL%1,               Roujian%Lu%NULL%1,               Peihua%Niu%NULL%1,               Faxian%Zhan%NULL%1,               Xuejun%Ma%NULL%1,               Dayan%Wang%NULL%1,               Wenbo%Xu%NULL%1,               Guizhen%Wu%NULL%1,               George F.%Gao%NULL%2,               Wenjie%Tan%NULL%2]'
$ws.Range("E3").Value = '[Xiaoxia%Lu%NULL%1,               Liqiong%Zhang%NULL%1,               Hui%Du%NULL%3,               Jingjing%Zhang%NULL%2,               Jingjing%Zhang%NULL%0,               Yuan Y.%Li%NULL%1,               Jingyu%Qu%NULL%1,               Wenxin%Zhang%NULL%1,               Youjie%Wang%NULL%1,               Shuangshuang%Bao%NULL%1,               Ying%Li%NULL%2,               Chuansha%Wu%NULL%1,               Hongxiu%Liu%NULL%1,               Di%Liu%NULL%0,               Di%Liu%NULL%0,               Jianbo%Shao%NULL%3,               Jianbo%Shao%NULL%0,               Jianbo%Shao%NULL%0,               Xuehua%Peng%NULL%1,               Yonghong%Yang%NULL%2,               Yonghong%Yang%NULL%0,               Zhisheng%Liu%NULL%3,               Zhisheng%Liu%NULL%0,               Zhisheng%Liu%NULL%0,               Yun%Xiang%NULL%1,               Furong%Zhang%NULL%1,               Rona M.%Silva%NULL%2,               Rona M.%Silva%NULL%0,               Kent E.%Pinkerton%NULL%1,               Kunling%Shen%NULL%2,               Kunling%Shen%NULL%0,               Han%Xiao%NULL%3,               Han%Xiao%NULL%0,               Shunqing%Xu%NULL%2,               Shunqing%Xu%NULL%0,               Gary W.K.%Wong%NULL%2,               Gary W.K.%Wong%NULL%0]'
$ws.Range("E5").Value = '[Esquivel%Sosa Leidelen%coreGivesNoEmail%1,              Mart\u00ednez-Fort\u00fan%Amador Maryla%coreGivesNoEmail%1,              \u00c1guila%Carbelo Madyaret%coreGivesNoEmail%1]'
$ws.Range("E7").Value = '[Benjamin%Lee%xref no email%1,        William V.%Raszka%xref no email%1]'
$ws.Range("E8").Value = '[Ruiyun%Li%NULL%1,               Sen%Pei%NULL%2,               Sen%Pei%NULL%0,               Bin%Chen%NULL%2,               Bin%Chen%NULL%0,               Yimeng%Song%NULL%2,               Yimeng%Song%NULL%0,               Tao%Zhang%NULL%3,               Tao%Zhang%NULL%0,               Wan%Yang%NULL%1,               Jeffrey%Shaman%NULL%2,               Jeffrey%Shaman%NULL%0]'
$ws.Range("E9").Value = '[Andrea%T. Cruz%coreGivesNoEmail%1,              Steven%L. Zeichner%coreGivesNoEmail%1]'
$ws.Range("E10").Value = '[Dan%Sun%NULL%3,               Hui%Li%NULL%0,               Xiao-Xia%Lu%NULL%1,               Han%Xiao%NULL%0,               Jie%Ren%NULL%1,               Fu-Rong%Zhang%792523496@qq.com%1,               Zhi-Sheng%Liu%liuzsc@126.com%1]'
$ws.Range("E11").Value = '[Hin%Chu%kyyuen@hku.hk%1,               Jie%Zhou%NULL%1,               Bosco Ho-Yin%Wong%NULL%1,               Cun%Li%NULL%1,               Jasper Fuk-Woo%Chan%NULL%0,               Zhong-Shan%Cheng%NULL%1,               Dong%Yang%NULL%1,               Dong%Wang%NULL%1,               Andrew Chak-Yiu%Lee%NULL%1,               Chuangen%Li%NULL%1,               Man-Lung%Yeung%NULL%1,               Jian-Piao%Cai%NULL%0,               Ivy Hau-Yee%Chan%NULL%1,               Wai-Kuen%Ho%NULL%1,               Kelvin Kai-Wang%To%NULL%2,               Bo-Jian%Zheng%NULL%1,               Yanfeng%Yao%NULL%1,               Chuan%Qin%NULL%2,               Kwok-Yung%Yuen%NULL%0]'
$ws.Range("E12").Value = '[William J.%Liu%liujun@ivdc.chinacdc.cn%1,               Min%Zhao%NULL%2,               Kefang%Liu%NULL%1,               Kun%Xu%NULL%1,               Gary%Wong%NULL%1,               Wenjie%Tan%NULL%0,               George F.%Gao%gaofu@chinacdc.cn%0]'
$ws.Range("E13").Value = '[Henry%Brandon Michael%coreGivesNoEmail%1,              Lippi%Giuseppe%coreGivesNoEmail%4,              Plebani%Mario%coreGivesNoEmail%1]'
